$wb = $excel.ActiveWorkbook

# --- Rename worksheet RM_01 -> BEN_FTR_001_Utilities ---
$wsUtil = $wb.Worksheets.Item("RM_01")
$wsUtil.Name = "BEN_FTR_001_Utilities"

# --- Common sheet: move selection from C17 to E3 ---
$wsCommon = $wb.Worksheets.Item("Common")
[void]$wsCommon.Range("E3").Select()

# --- BEN_FTR_001_Utilities sheet: move selection from B8 to C10 ---
[void]$wsUtil.Range("C10").Select()

# --- Credentials sheet: add a new Recruitment user row + hyperlink ---
$wsCred = $wb.Worksheets.Item("Credentials")

$wsCred.Range("B4").Value = "KMizell@SEMPRANRGU "
$wsCred.Range("C4").Value = "ProKarma2018!"
$wsCred.Range("A4").Value = "Recruitment2"
$wsCred.Range("E4").Value = "This is Ken's user for Recruitment"
$wsCred.Range("D4").Value = "Yes"

[void]$wsCred.Hyperlinks.Add($wsCred.Range("B4"), "mailto:KMizell@SEMPRANRGU")

# Widen column C a bit to fit the new password value (best-effort match)
$wsCred.Columns.Item(3).ColumnWidth = 13.6

# Return focus/selection to the Credentials sheet, cell A1
[void]$wsCred.Range("A1").Select()
